$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 236-237, pushing the existing 236-237 (week of
# 44432) down to 238-239, so the new weekly data lands on top.
$ws.Rows("236:237").Insert()

# New row 236: Apio, Americana (o), Primera - week of 44656
$ws.Cells.Item(236, 1).Value = 11
$ws.Cells.Item(236, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(236, 3).Value = "Bíobío"
$ws.Cells.Item(236, 4).Value = 44656
$ws.Cells.Item(236, 5).Value = 8
$ws.Cells.Item(236, 6).Value = 100112017
$ws.Cells.Item(236, 7).Value = "Apio"
$ws.Cells.Item(236, 8).Value = "Americana (o)"
$ws.Cells.Item(236, 9).Value = "Primera"
$ws.Cells.Item(236, 10).Value = 250
$ws.Cells.Item(236, 11).Value = 8000
$ws.Cells.Item(236, 12).Value = 8500
$ws.Cells.Item(236, 13).Value = 8260
$ws.Cells.Item(236, 14).Value = "`$/docena de matas"
$ws.Cells.Item(236, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(236, 16).Value = 1377
$ws.Cells.Item(236, 17).Value = 6
$ws.Cells.Item(236, 18).Value = "Hortaliza"

# New row 237: Apio, Americana (o), Segunda - week of 44656
$ws.Cells.Item(237, 1).Value = 11
$ws.Cells.Item(237, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(237, 3).Value = "Bíobío"
$ws.Cells.Item(237, 4).Value = 44656
$ws.Cells.Item(237, 5).Value = 8
$ws.Cells.Item(237, 6).Value = 100112017
$ws.Cells.Item(237, 7).Value = "Apio"
$ws.Cells.Item(237, 8).Value = "Americana (o)"
$ws.Cells.Item(237, 9).Value = "Segunda"
$ws.Cells.Item(237, 10).Value = 220
$ws.Cells.Item(237, 11).Value = 6000
$ws.Cells.Item(237, 12).Value = 6500
$ws.Cells.Item(237, 13).Value = 6273
$ws.Cells.Item(237, 14).Value = "`$/docena de matas"
$ws.Cells.Item(237, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(237, 16).Value = 1046
$ws.Cells.Item(237, 17).Value = 6
$ws.Cells.Item(237, 18).Value = "Hortaliza"
